# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# sheet with refreshed figures from the latest GitHub Actions run.
# Numeric-looking "Price" values are prefixed with a leading apostrophe
# so Excel keeps storing them as text (matching the source data, which
# uses dotted thousand separators like "42.711.81" that are not valid
# numbers) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.711.81'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.262.07'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''249.24'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = '''0.636'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('D7').Value = '''78.39'
$ws.Range('E7').Value = '  +6.81%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = '''0.642'
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('D10').Value = '''40.87'
$ws.Range('E10').Value = '  +3.76%  '
$ws.Range('D11').Value = '''0.0958'
$ws.Range('E11').Value = '  -2.38%  '
$ws.Range('D12').Value = '''7.30'
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('D14').Value = '2.601.85'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '''14.98'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = '''0.864'
$ws.Range('E16').Value = '  -4.10%  '
$ws.Range('D17').Value = '2.277.55'
$ws.Range('E17').Value = '  -1.20%  '
$ws.Range('D18').Value = '42.468.90'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').Value = '0.0₃0988'
$ws.Range('D20').Value = '''6.17'
$ws.Range('E20').Value = '  -3.92%  '
$ws.Range('D21').Value = '''71.79'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').Value = '''232.63'
$ws.Range('E22').Value = '  -2.35%  '
$ws.Range('D23').Value = '''2.15'
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').Value = '''3.79'
$ws.Range('E24').Value = '  -3.57%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '''11.30'
$ws.Range('E26').Value = '  -4.83%  '
$ws.Range('D27').Value = '''2.31'
$ws.Range('E27').Value = '  -5.51%  '
$ws.Range('D28').Value = '''2.17'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').Value = '''170.19'
$ws.Range('E29').Value = '  +1.16%  '
$ws.Range('D30').Value = '''6.65'
$ws.Range('E30').Value = '  +4.44%  '
$ws.Range('D31').Value = '''20.72'
$ws.Range('E31').Value = '  -2.81%  '
$ws.Range('D32').Value = '''0.0846'
$ws.Range('E32').Value = '  +3.90%  '
$ws.Range('D33').Value = '''0.121'
$ws.Range('E33').Value = '  -5.83%  '
$ws.Range('D34').Value = '''30.54'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('E35').Value = '  -0.42%  '
$ws.Range('D36').Value = '''4.52'
$ws.Range('E36').Value = '  -4.05%  '
$ws.Range('D37').Value = '''4.72'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').Value = '''0.0301'
$ws.Range('E38').Value = '  -3.61%  '
$ws.Range('D39').Value = '''13.31'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').Value = '''2.23'
$ws.Range('E40').Value = '  -5.62%  '
$ws.Range('D41').Value = '''5.96'
$ws.Range('E41').Value = '  -2.90%  '
$ws.Range('D42').Value = '''113.98'
$ws.Range('E42').Value = '  +16.67%  '
$ws.Range('D43').Value = '''0.207'
$ws.Range('E43').Value = '  -3.03%  '
$ws.Range('D44').Value = '''60.95'
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D45').Value = '''8.85'
$ws.Range('E45').Value = '  -4.87%  '
$ws.Range('E46').Value = '  -2.83%  '
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('D48').Value = '''4.48'
$ws.Range('E48').Value = '  -8.93%  '
$ws.Range('D49').Value = '''1.14'
$ws.Range('E49').Value = '  -4.69%  '
$ws.Range('D50').Value = '''1.16'
$ws.Range('E50').Value = '  -3.26%  '
$ws.Range('D51').Value = '''4.20'
$ws.Range('E51').Value = '  -2.39%  '
